$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^\$?\s*-?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

$data = @{
    833 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    834 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    835 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    836 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    837 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    838 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    839 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM' }
    840 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM' }
    841 = @{ 'A'='21CRB01268'; 'B'='Hemmeter'; 'C'='POSSESSION DRUG PARAPHERNALIA'; 'D'='2925.14(C)'; 'E'='M4'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    842 = @{ 'A'='21CRB01268'; 'B'='Hemmeter'; 'C'='POSSESSION DRUG PARAPHERNALIA'; 'D'='2925.14(C)'; 'E'='M4'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    843 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='Guilty'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    844 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='Guilty'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    845 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='Guilty'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    846 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='Guilty'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    847 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    848 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    849 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    850 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    851 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    852 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    853 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    854 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    855 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='DUS UCM'; 'D'='4510.111'; 'E'='UCM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
    856 = @{ 'A'='21TRD09386'; 'B'='Hemmeter'; 'C'='TAIL LIGHTS-REAR LICENSE PLATE'; 'D'='4513.05'; 'E'='MM'; 'F'='No Contest'; 'G'='Guilty'; 'H'='$ 0'; 'I'='$ 0' }
}

foreach ($r in $data.Keys) {
    $rowData = $data[$r]
    foreach ($c in $rowData.Keys) {
        $colIndex = switch ($c) {
            'A' { 1 }
            'B' { 2 }
            'C' { 3 }
            'D' { 4 }
            'E' { 5 }
            'F' { 6 }
            'G' { 7 }
            'H' { 8 }
            'I' { 9 }
        }
        Set-CellText ([int]$r) $colIndex $rowData[$c]
    }
}
